$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $text) {
    $row = $t.Rows.Item($rowIndex)
    $cell = $row.Cells.Item(1)
    $cell.Range.Text = $text
}

# Single-value summary cells near the top of the table.
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "306"
Set-CellText 5 "0.00002"
Set-CellText 6 "0.00072"
Set-CellText 7 "0.00017"
Set-CellText 9 "0.00036"
Set-CellText 10 "0.00040"
Set-CellText 11 "0.00042"
Set-CellText 12 "0.06232"

# Rows near the bottom that previously held tab-separated multi-column
# data get collapsed down to a single value each.
Set-CellText 44 "99.86"
Set-CellText 45 "0.06"
Set-CellText 46 "45"
